$wb = $excel.ActiveWorkbook

# Fix the escaped quotes in the "ratio de Sharpe" question text.
# This text is duplicated (by value) across pregunta!B2 and respuesta!F2:F4.
$fixedText = '¿Qué significa el término "ratio de Sharpe" y cómo se utiliza en la evaluación de carteras de inversión?'

$wsPregunta = $wb.Worksheets.Item("pregunta")
$wsPregunta.Range("B2").Value = $fixedText

$wsRespuesta = $wb.Worksheets.Item("respuesta")
$wsRespuesta.Range("F2").Value = $fixedText
$wsRespuesta.Range("F3").Value = $fixedText
$wsRespuesta.Range("F4").Value = $fixedText
